$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet stores one data row per market observation. The edit inserts
# three brand-new observation rows right after the existing header+data rows
# at positions 1070-1072 (pushing everything that used to live at rows
# 1070-1142 down by three rows, to 1073-1145), and fills the three newly
# freed rows (1070-1072) with new data.
#
# To avoid clobbering source data before it is copied, we must copy rows in
# DESCENDING order: target row N (N going from 1145 down to 1073) gets the
# content that used to be in row (N-3). Because N-3 < N always, and we work
# from the highest N down to the lowest, we never overwrite a row before we
# have read it.
# ---------------------------------------------------------------------------

$firstOldRow = 1070
$lastOldRow  = 1142
$shift       = 3
$lastCol     = 18   # columns A..R

for ($n = $lastOldRow + $shift; $n -ge $firstOldRow + $shift; $n--) {
    $src = $n - $shift
    for ($c = 1; $c -le $lastCol; $c++) {
        $v = $ws.Cells.Item($src, $c).Value2
        $ws.Cells.Item($n, $c).Value2 = $v
    }
    # Column D (dates) carries a custom number format (style index 2 in the
    # original file). Rows beyond the old last row (1143-1145) are brand new
    # cells with no formatting yet, so copy the date format explicitly.
    if ($n -gt $lastOldRow) {
        $ws.Cells.Item($n, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
    }
}

# ---------------------------------------------------------------------------
# Now fill in the three newly freed rows (1070, 1071, 1072) with their new
# data values.
# ---------------------------------------------------------------------------

# Row 1070
$ws.Cells.Item(1070, 1).Value2  = 3
$ws.Cells.Item(1070, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(1070, 3).Value2  = "Coquimbo"
$ws.Cells.Item(1070, 4).Value2  = 44714
$ws.Cells.Item(1070, 5).Value2  = 5
$ws.Cells.Item(1070, 6).Value2  = 100112020
$ws.Cells.Item(1070, 7).Value2  = "Tomate"
$ws.Cells.Item(1070, 8).Value2  = "Larga vida"
$ws.Cells.Item(1070, 9).Value2  = "Primera"
$ws.Cells.Item(1070, 10).Value2 = 540
$ws.Cells.Item(1070, 11).Value2 = 13500
$ws.Cells.Item(1070, 12).Value2 = 14000
$ws.Cells.Item(1070, 13).Value2 = 13759
$ws.Cells.Item(1070, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(1070, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(1070, 16).Value2 = 764
$ws.Cells.Item(1070, 17).Value2 = 18
$ws.Cells.Item(1070, 18).Value2 = "Hortaliza"

# Row 1071
$ws.Cells.Item(1071, 1).Value2  = 3
$ws.Cells.Item(1071, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(1071, 3).Value2  = "Coquimbo"
$ws.Cells.Item(1071, 4).Value2  = 44714
$ws.Cells.Item(1071, 5).Value2  = 5
$ws.Cells.Item(1071, 6).Value2  = 100112020
$ws.Cells.Item(1071, 7).Value2  = "Tomate"
$ws.Cells.Item(1071, 8).Value2  = "Larga vida"
$ws.Cells.Item(1071, 9).Value2  = "Primera"
$ws.Cells.Item(1071, 10).Value2 = 770
$ws.Cells.Item(1071, 11).Value2 = 6500
$ws.Cells.Item(1071, 12).Value2 = 7000
$ws.Cells.Item(1071, 13).Value2 = 6753
$ws.Cells.Item(1071, 14).Value2 = "$/caja 12 kilos"
$ws.Cells.Item(1071, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(1071, 16).Value2 = 563
$ws.Cells.Item(1071, 17).Value2 = 12
$ws.Cells.Item(1071, 18).Value2 = "Hortaliza"

# Row 1072
$ws.Cells.Item(1072, 1).Value2  = 3
$ws.Cells.Item(1072, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(1072, 3).Value2  = "Coquimbo"
$ws.Cells.Item(1072, 4).Value2  = 44714
$ws.Cells.Item(1072, 5).Value2  = 5
$ws.Cells.Item(1072, 6).Value2  = 100112020
$ws.Cells.Item(1072, 7).Value2  = "Tomate"
$ws.Cells.Item(1072, 8).Value2  = "Larga vida"
$ws.Cells.Item(1072, 9).Value2  = "Segunda"
$ws.Cells.Item(1072, 10).Value2 = 560
$ws.Cells.Item(1072, 11).Value2 = 11500
$ws.Cells.Item(1072, 12).Value2 = 12000
$ws.Cells.Item(1072, 13).Value2 = 11750
$ws.Cells.Item(1072, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(1072, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(1072, 16).Value2 = 653
$ws.Cells.Item(1072, 17).Value2 = 18
$ws.Cells.Item(1072, 18).Value2 = "Hortaliza"
